$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Updated salinity data (column U) -------------------------------------
# Rows 2-8 and 12-21 had salinity re-measured/re-run from 35 -> 32.
# (Rows 9-11 and 22-23 already hold distinct values and are left untouched.)
$salinityRows = @(2,3,4,5,6,7,8,12,13,14,15,16,17,18,19,20,21)
foreach ($r in $salinityRows) {
    $ws.Range("U$r").Value = 32
}

# --- View state -------------------------------------------------------------
# Re-ran script: scrolled/selected a different area of the sheet before saving.
$ws.Activate()
$excel.ActiveWindow.ScrollRow = 7
$excel.ActiveWindow.ScrollColumn = 6
$ws.Range("U27:V27").Select()
